# Soldering station deck edit:
#  - remove the "Libraries and Schematics" slide (old slide 4 / sldId 259)
#  - hide the "Content" agenda slide
#  - reposition the "Introduction" title and drop its empty content placeholder
#  - reposition the "Final product" title and drop its "Demo" content placeholder

$p = $ppt.ActivePresentation

# --- Delete the "Libraries and Schematics" slide (position 4) ---------------
$p.Slides.Item(4).Delete()

# --- Hide the "Content" agenda slide (position 2) ---------------------------
$content = $p.Slides.Item(2)
$content.SlideShowTransition.Hidden = -1

# --- "Introduction" slide (position 3): move title, drop empty placeholder --
$intro = $p.Slides.Item(3)
$introTitle = $intro.Shapes.Item(1)
$introTitle.Left = 166.67984251968505
$introTitle.Top = 227.58937007874016
$intro.Shapes.Item(2).Delete()

# --- "Final product" slide (now last, position 6): move title, drop "Demo" --
$final = $p.Slides.Item($p.Slides.Count)
$finalTitle = $final.Shapes.Item(1)
$finalTitle.Left = 166.67984251968505
$finalTitle.Top = 227.58937007874016
$finalDemo = $final.Shapes.Item(2)
$finalDemo.TextFrame.TextRange.Text = ""
$finalDemo.Delete()
